$d = $word.ActiveDocument

# Locate the paragraph containing "git branch -d feature-1" (the one right
# after "To delete a feature branch"), using Find so we do not depend on a
# hardcoded paragraph index.
$findRange = $d.Content
$ok = $findRange.Find.Execute("git branch -d feature-1", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find target paragraph 'git branch -d feature-1'"
}

# Figure out which paragraph (1-based index into $d.Paragraphs) holds the hit.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs($i).Range
    if ($pr.Start -le $findRange.Start -and $findRange.Start -lt $pr.End) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq 0) {
    throw "Could not resolve paragraph index for the match"
}

$origPara = $d.Paragraphs($targetIndex)
$origRange = $origPara.Range

# Insert two new (initially empty) paragraphs right before the original one.
# Word duplicates the paragraph-mark formatting of the following paragraph
# into each newly created paragraph; the original paragraph (and its own
# paragraph-mark formatting) ends up shifted two slots later.
$origRange.InsertParagraphBefore()
$origRange.InsertParagraphBefore()

$p1 = $d.Paragraphs($targetIndex)
$p2 = $d.Paragraphs($targetIndex + 1)
$p3 = $d.Paragraphs($targetIndex + 2)

# --- Paragraph 1: "git branch -d feature-1" (same text/format as before,
#     now living in its own paragraph with a "normal" 14pt paragraph mark).
$p1.Range.Font.Size = 14
$p1.Range.InsertBefore("git branch -d feature-1")

# --- Paragraph 2: new explanatory line, same (non-bold) formatting.
$p2.Range.Font.Size = 14
$p2.Range.InsertBefore("To delete a branch directly from github")

# --- Paragraph 3: reuse the original paragraph (keeps its original
#     paragraph-mark formatting) but replace its run text with the new bold
#     "git push origin -d Lesson7-Porject1-InteractiveQuiz" command.
$p3Range = $p3.Range
$textOnly = $d.Range($p3Range.Start, $p3Range.End - 1)
$textOnly.Text = "git push origin -d "
$textOnly.Font.Bold = $true

$p3Range = $p3.Range
$insertAt = $p3Range.End - 1
$tail = $d.Range($insertAt, $insertAt)
$tail.InsertAfter("Lesson7-Porject1-InteractiveQuiz")
$tailFormatted = $d.Range($insertAt, $insertAt + "Lesson7-Porject1-InteractiveQuiz".Length)
$tailFormatted.Font.Bold = $true

Write-Host "p1: [$($p1.Range.Text)]"
Write-Host "p2: [$($p2.Range.Text)]"
Write-Host "p3: [$($p3.Range.Text)]"
